$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title heading: "Some Class Name" -> "Librarian"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Some Class Name", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Librarian", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. CRC table header cell: "Class Name" -> "Librarian"
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)
$cell = $tbl.Cell(1, 1)
$cell.Range.Find.Execute("Class Name", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Librarian", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. The description paragraph used to contain the run sequence
#       "So" + bookmarkStart(_GoBack) + bookmarkEnd(_GoBack) + "me paragraph about this class"
#    which together read "Some paragraph about this class". Word has since
#    re-typed over that text, so the split runs/bookmark collapse back into a
#    single plain run and the "_GoBack" bookmark moves to mark the most
#    recently edited location - right after the class name we just typed in
#    the CRC table cell above.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$p = $goBack.Start

# Temporary bookmarks fence off exactly the two runs ("So" / "me paragraph
# about this class") that need to be recombined, so the text reflow below
# only touches that span and leaves every other run in the paragraph alone.
$fenceStart = $d.Range($p - 2, $p - 2)
$d.Bookmarks.Add("ZZZ_FENCE_START", $fenceStart) | Out-Null
$fenceEnd = $d.Range($p + 29, $p + 29)
$d.Bookmarks.Add("ZZZ_FENCE_END", $fenceEnd) | Out-Null

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$span = $d.Range($p - 2, $p + 29)
$span.Text = "Some paragraph about this class~~RETYPED~~"

$startFence = $d.Bookmarks("ZZZ_FENCE_START")
$endFence = $d.Bookmarks("ZZZ_FENCE_END")
$merged = $d.Range($startFence.Start, $endFence.End)
$merged.Text = "Some paragraph about this class"

$d.Bookmarks("ZZZ_FENCE_START").Delete()
$d.Bookmarks("ZZZ_FENCE_END").Delete()

# ---------------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark at the last-edited spot: right after
#    the "Librarian" text that now fills the CRC table's class-name cell.
# ---------------------------------------------------------------------------
$cell = $d.Tables(1).Cell(1, 1)
$dest = $cell.Range
$dest.Collapse(0)
$dest.MoveEnd(1, -1) | Out-Null
$dest.Collapse(0)
$d.Bookmarks.Add("_GoBack", $dest) | Out-Null
